$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Wipe the old scratch row (A2:M2 values) and the two empty bordered
#     cells (A1/A3) that don't belong in the new layout. Re-seating the
#     style to the plain default (instead of just ClearFormats) is what
#     actually drops these cells back out of the sheet once they're empty. ---
$defaultStyle = $ws.Range("B2").Style
$ws.Range("A2").Style = $defaultStyle
$ws.Range("A3").Style = $defaultStyle
$ws.Range("A2:M2").ClearContents()
$ws.Range("A3").ClearContents()

# --- Header row (row 1): "Unnamed: 0" / "2019" / "Unnamed: 1", all using the
#     bordered/bold/centered style that already lives on A1. ---
# "2019" reads as a number, so the cell has to be pre-formatted as Text,
# otherwise typing it in auto-converts it to a numeric value.
$ws.Range("B1").NumberFormat = "@"

$ws.Range("A1").Value = "Unnamed: 0"
$ws.Range("B1").Value = "2019"
$ws.Range("C1").Value = "Unnamed: 1"

# Copy A1's formatting (font/border/alignment) onto B1:C1 so the whole
# header row matches.
$ws.Range("A1").Copy()
$ws.Range("B1:C1").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# --- Data column (column B): rows 2-4 then 6-14 (row 5 stays blank). ---
$values = @(186, 101, 85, 57, 31, 26, 144, 97, 47, 35, 23, 12)
$rows = @(2, 3, 4, 6, 7, 8, 9, 10, 11, 12, 13, 14)
for ($i = 0; $i -lt $values.Length; $i++) {
    $ws.Cells.Item($rows[$i], 2).Value = $values[$i]
}
